$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'59.364.40"
$ws.Range("E2").Value = "  +0.73%  "

$ws.Range("D3").Value = "'2.527.37"
$ws.Range("E3").Value = "  +1.30%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").Value = "'535.76"
$ws.Range("E5").Value = "  -0.22%  "

$ws.Range("D6").Value = "'140.18"

$ws.Range("E7").Value = "  +0.16%  "

$ws.Range("D8").Value = "'0.563"
$ws.Range("E8").Value = "  -1.91%  "

$ws.Range("D9").Value = "'2.532.94"
$ws.Range("E9").Value = "  +0.51%  "

$ws.Range("E10").Value = "  -0.29%  "

$ws.Range("E11").Value = "  +1.74%  "

$ws.Range("E12").Value = "  -0.36%  "

$ws.Range("D13").Value = "'0.355"
$ws.Range("E13").Value = "  +0.82%  "

$ws.Range("D14").Value = "'2.975.34"
$ws.Range("E14").Value = "  +1.39%  "

$ws.Range("D15").Value = "'23.15"
$ws.Range("E15").Value = "  -2.46%  "

$ws.Range("D16").Value = "'59.301.97"
$ws.Range("E16").Value = "  +0.79%  "

$ws.Range("E17").Value = "  +1.40%  "

$ws.Range("D18").Value = "'2.548.35"
$ws.Range("E18").Value = "  +1.44%  "

$ws.Range("D19").Value = "'10.97"
$ws.Range("E19").Value = "  -2.70%  "

$ws.Range("E20").Value = "  -1.16%  "

$ws.Range("D21").Value = "'321.35"
$ws.Range("E21").Value = "  -0.25%  "

$ws.Range("E22").Value = "  +0.29%  "

$ws.Range("D23").Value = "'5.83"
$ws.Range("E23").Value = "  +1.48%  "

$ws.Range("D24").Value = "'62.18"
$ws.Range("E24").Value = "  +0.68%  "

$ws.Range("E25").Value = "  -4.02%  "

$ws.Range("E26").Value = "  +2.47%  "

$ws.Range("D27").Value = "'0.997"
$ws.Range("E27").Value = "  +0.32%  "

$ws.Range("D28").Value = "'7.83"

$ws.Range("D29").Value = "'6.74"
$ws.Range("E29").Value = "  -0.25%  "

$ws.Range("E30").Value = "  -0.76%  "

$ws.Range("E31").Value = "  +0.37%  "

$ws.Range("D32").Value = "'161.23"
$ws.Range("E32").Value = "  +1.22%  "

$ws.Range("E33").Value = "  +0.28%  "

$ws.Range("E34").Value = "  -5.72%  "

$ws.Range("D35").Value = "'1.46"
$ws.Range("E35").Value = "  -0.27%  "

$ws.Range("D36").Value = "'18.52"
$ws.Range("E36").Value = "  -0.08%  "

$ws.Range("E37").Value = "  -3.15%  "

$ws.Range("D38").Value = "'1.59"
$ws.Range("E38").Value = "  -2.02%  "

$ws.Range("D39").Value = "'37.07"
$ws.Range("E39").Value = "  +0.71%  "

$ws.Range("E40").Value = "  -0.36%  "

$ws.Range("B41").Value = "SuiNetwork"
$ws.Range("C41").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D41").Value = "'0.808"
$ws.Range("E41").Value = "  -1.60%  "

$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D42").Value = "'5.30"
$ws.Range("E42").Value = "  -7.10%  "

$ws.Range("B43").Value = "Bittensor"
$ws.Range("C43").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D43").Value = "'286.57"
$ws.Range("E43").Value = "  -5.87%  "

$ws.Range("E44").Value = "  +0.52%  "

$ws.Range("E45").Value = "  +0.97%  "

$ws.Range("E46").Value = "  +0.69%  "

$ws.Range("D47").Value = "'124.27"
$ws.Range("E47").Value = "  -0.92%  "

$ws.Range("E48").Value = "  -0.29%  "

$ws.Range("D49").Value = "'18.58"
$ws.Range("E49").Value = "  +0.06%  "

$ws.Range("E50").Value = "  -1.46%  "

$ws.Range("E51").Value = "  -2.18%  "
